$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Partida 15): Atlantico x Galo Futsal result entered -> 6x2, Finalizado
$ws.Range("E16").Value = "6x2"
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = "Finalizado"

# Row 17 (Partida 16): Real Carira x Vasquinho result entered -> 2x4, Finalizado
$ws.Range("E17").Value = "2x4"
$ws.Range("G17").Value = 2
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = "Finalizado"

# Update the view state: scrolled so row 7 is the top row, and N17 is selected
$ws.Range("N17").Select()
$excel.ActiveWindow.ScrollRow = 7
